# Se actualiza los puntos realizados por los trabajadores
$wb = $excel.ActiveWorkbook

$wsCel     = $wb.Worksheets.Item("Mes actual cel")
$wsEmp     = $wb.Worksheets.Item("Mes actual emp")
$wsHistCel = $wb.Worksheets.Item("Historico CEL")

# --- "Mes actual cel": puntos del mes actual por trabajador de cuadrilla ---
$wsCel.Range("C2").Value = 58.94
$wsCel.Range("C3").Value = 62.67
$wsCel.Range("C4").Value = 58.94

# --- "Mes actual emp": puntos del mes actual del empleado ---
$wsEmp.Range("C2").Value = 84.1688
$wsEmp.Range("C2").NumberFormat = "0.00"

# --- vistas/selecciones/zoom tal y como quedaron grabadas por Excel ---
$wsEmp.Activate()
$wsEmp.Range("E18").Select()
$excel.ActiveWindow.Zoom = 220

$wsHistCel.Activate()
$wsHistCel.Range("J23").Select()
$excel.ActiveWindow.Zoom = 205

$wsCel.Activate()
$wsCel.Range("I10").Select()
$excel.ActiveWindow.Zoom = 250
